$d = $word.ActiveDocument

# 1) "Reconhecimento e reconhecimento" / " de marca limitados: ..."
$d.Content.Find.Execute("Reconhecimento e reconhecimento", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reconhecimento e conscientização da marca limitados", 2)
$d.Content.Find.Execute(" de marca limitados: Alcançar visibilidade nesses novos mercados é um obstáculo primário, exigindo esforços de marketing robustos para construir a presença da marca Adatum desde o início.", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": alcançar visibilidade nesses novos mercados é um grande obstáculo, exigindo esforços robustos de marketing para construir a presença de marca da Adatum do zero.", 2)

# 2) "Concorrência" / " intensa: ..."
$d.Content.Find.Execute("Concorrência", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Concorrência intensa", 2)
$d.Content.Find.Execute(" intensa: O setor de serviços em nuvem no Canadá é ferozmente competitivo, com vários players.", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": o setor de serviços em nuvem no Canadá é ferozmente competitivo, com vários envolvidos.", 2)

# 3) "Preferências e expectativas diversificadas" / " dos clientes: ..."
$d.Content.Find.Execute("Preferências e expectativas diversificadas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Preferências e expectativas diversas do cliente", 2)
$d.Content.Find.Execute(" dos clientes: Adaptar os produtos e o marketing para se alinhar às variadas demandas desses mercados é crucial para ressoar com as empresas e consumidores locais.", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": adaptar produtos e marketing para se alinhar às diversas demandas desses mercados é crucial para identificação com as empresas e os consumidores locais.", 2)

# 4) "Desafios" / " regulatórios e de conformidade: ..."
$d.Content.Find.Execute("Desafios", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Desafios de regulamentação e conformidade", 2)
$d.Content.Find.Execute(" regulatórios e de conformidade: A Adatum enfrenta a complexa tarefa de navegar pelas distintas regulamentações de privacidade, segurança e operação de dados da região, exigindo esforços diligentes de conformidade.", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": a Adatum enfrenta a complexa tarefa de navegar pelos distintos regulamentos operacionais e de privacidade de dados da região, o que demanda esforços diligentes de conformidade.", 2)

# 5) "Complexidades operacionais e logísticas" stays the same; only the trailing run changes
$d.Content.Find.Execute(": O estabelecimento de operações eficientes e inter-regionais apresenta desafios logísticos, especialmente na manutenção de altos níveis de serviço e no gerenciamento de data centers em localizações geográficas.", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": estabelecer operações inter-regionais eficientes apresenta desafios logísticos, especialmente na manutenção de altos níveis de serviço e no gerenciamento de data centers em localizações geográficas.", 2)
